$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values (column B) that changed text
$ws.Range("B3").Value = "C:/Users/mt5285/ptx_data/settings/test.xlsx"
$ws.Range("B4").Value = "visualize_only"
$ws.Range("B6").Value = "C:/Users/mt5285/ptx_data/settings/porsche/"

# Add new row 9: path_visualization / result path
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A9").Value = "path_visualization"

$ws.Range("B9").Value = "C:/Users/mt5285/ptx_data/results/20211220_093208_test/"
